$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7100
$ws.Range("I51").Value = 7075
$ws.Range("K51").Value = 7075
$ws.Range("M51").Value = -6591

$ws.Range("H64").Value = 7660.6
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 8575.75
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 8575.75
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -9071.75

$ws.Range("H67").Value = 7660.6
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 8575.75
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 8575.75
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -10291.75

$ws.Range("H76").Value = 3665.5
$ws.Range("I76").Value = 3748.5
$ws.Range("J76").Value = 3624
$ws.Range("K76").Value = 3748.5
$ws.Range("L76").Value = 3624
$ws.Range("M76").Value = -3433.5
$ws.Range("N76").Value = -4254

$ws.Range("H79").Value = 3665.5
$ws.Range("I79").Value = 3748.5
$ws.Range("J79").Value = 3624
$ws.Range("K79").Value = 3748.5
$ws.Range("L79").Value = 3624
$ws.Range("M79").Value = -2656.5
$ws.Range("N79").Value = -5808

$ws.Range("H98").Value = 574.7143
$ws.Range("I98").Value = 542
$ws.Range("K98").Value = 542
$ws.Range("M98").Value = 956

$ws.Range("H122").Value = 574.7143
$ws.Range("I122").Value = 542
$ws.Range("K122").Value = 1626
$ws.Range("M122").Value = 824

$ws.Range("H125").Value = 1094
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H137").Value = 3903.7083
$ws.Range("I137").Value = 3110.3076
$ws.Range("J137").Value = 4841.364
$ws.Range("K137").Value = 9330.9228
$ws.Range("L137").Value = 14524.092
$ws.Range("M137").Value = -6780.9228
$ws.Range("N137").Value = -19624.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1436.6666
$ws.Range("I2").Value = 991.5714
$ws.Range("K2").Value = 991.5714
$ws.Range("M2").Value = -878.5714

$ws.Range("H116").Value = 1436.6666
$ws.Range("I116").Value = 991.5714
$ws.Range("K116").Value = 991.5714
$ws.Range("M116").Value = 1302.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1436.6666
$ws.Range("I3").Value = 991.5714
$ws.Range("K3").Value = 991.5714
$ws.Range("M3").Value = -877.5714

$ws.Range("H12").Value = 580
$ws.Range("I12").Value = 610
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 610
$ws.Range("L12").Value = 400
$ws.Range("M12").Value = -442
$ws.Range("N12").Value = -736

$ws.Range("H105").Value = 2425
$ws.Range("I105").Value = 1900
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 1900
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -153
$ws.Range("N105").Value = -7494

$ws.Range("H107").Value = 5011.3335
$ws.Range("I107").Value = 3642
$ws.Range("K107").Value = 3642
$ws.Range("M107").Value = -1722

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4735.871
$ws.Range("I31").Value = 946.4375
$ws.Range("J31").Value = 8777.933999999999
$ws.Range("K31").Value = 946.4375
$ws.Range("L31").Value = 8777.933999999999
$ws.Range("M31").Value = -651.4375
$ws.Range("N31").Value = -9367.933999999999

$ws.Range("H34").Value = 4735.871
$ws.Range("I34").Value = 946.4375
$ws.Range("J34").Value = 8777.933999999999
$ws.Range("K34").Value = 946.4375
$ws.Range("L34").Value = 8777.933999999999
$ws.Range("M34").Value = -744.4375
$ws.Range("N34").Value = -9181.933999999999

$ws.Range("H105").Value = 2043.8667
$ws.Range("I105").Value = 1567.5555
$ws.Range("J105").Value = 2758.3333
$ws.Range("K105").Value = 1567.5555
$ws.Range("L105").Value = 2758.3333
$ws.Range("M105").Value = 179.4445000000001
$ws.Range("N105").Value = -6252.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1000000
$ws.Range("I11").Value = 1000000
$ws.Range("K11").Value = 3000000
$ws.Range("M11").Value = -2999860

$ws.Range("H68").Value = 2090.3684
$ws.Range("I68").Value = 1366.6666
$ws.Range("J68").Value = 2226.0625
$ws.Range("K68").Value = 4099.9998
$ws.Range("L68").Value = 6678.1875
$ws.Range("M68").Value = -3288.9998
$ws.Range("N68").Value = -8300.1875

$ws.Range("H71").Value = 2090.3684
$ws.Range("I71").Value = 1366.6666
$ws.Range("J71").Value = 2226.0625
$ws.Range("K71").Value = 12299.9994
$ws.Range("L71").Value = 20034.5625
$ws.Range("M71").Value = -8243.999400000001
$ws.Range("N71").Value = -28146.5625

$ws.Range("H103").Value = 561
$ws.Range("I103").Value = 349.4
$ws.Range("J103").Value = 1090
$ws.Range("K103").Value = 1048.2
$ws.Range("L103").Value = 3270
$ws.Range("M103").Value = -169.1999999999998
$ws.Range("N103").Value = -5028

$ws.Range("H107").Value = 324.57144
$ws.Range("I107").Value = 330.5
$ws.Range("J107").Value = 289
$ws.Range("K107").Value = 991.5
$ws.Range("L107").Value = 867
$ws.Range("M107").Value = 928.5
$ws.Range("N107").Value = -4707

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 268.95456
$ws.Range("I2").Value = 135.4
$ws.Range("K2").Value = 135.4
$ws.Range("M2").Value = -22.40000000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3205.111
$ws.Range("I122").Value = 3318.25
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 9954.75
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -7504.75
$ws.Range("N122").Value = -11800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3277.5
$ws.Range("I81").Value = 1000
$ws.Range("K81").Value = 2000
$ws.Range("M81").Value = -939

$ws.Range("H84").Value = 3277.5
$ws.Range("I84").Value = 1000
$ws.Range("K84").Value = 10000
$ws.Range("M84").Value = -4696

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H136").Value = 4483
$ws.Range("I136").Value = 4224.5
$ws.Range("K136").Value = 12673.5
$ws.Range("M136").Value = -10123.5
